$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) with same formatting as the existing header (H1: bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells row 2
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

# New data cells row 3
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 5
